$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.475208
$ws.Range("N2").Value = 4.425624
$ws.Range("O2").Value = 0.8210007041987012
$ws.Range("P2").Value = 0.8210007041987013
$ws.Range("Q2").Value = 54.500621819448
$ws.Range("R2").Value = 490.505596375032
$ws.Range("S2").Value = 0.6838086110431338
$ws.Range("T2").Value = 0.683808611043134

$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("M3").Value = 0.3216333333333333
$ws.Range("N3").Value = 0.9649
$ws.Range("O3").Value = 0.1789992958012987
$ws.Range("P3").Value = 0.1789992958012987
$ws.Range("Q3").Value = 11.88253904841111
$ws.Range("R3").Value = 106.9428514357
$ws.Range("S3").Value = 0.1490878865433484
$ws.Range("T3").Value = 0.1490878865433485

$ws.Range("I4").Value = 0.07608399754092349
$ws.Range("J4").Value = 0.07608399754092349
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.475208
$ws.Range("N4").Value = 4.425624
$ws.Range("O4").Value = 0.8210007041987012
$ws.Range("P4").Value = 0.8210007041987013
$ws.Range("Q4").Value = 4.978559987352
$ws.Range("R4").Value = 44.807039886168
$ws.Range("S4").Value = 0.06246501555935043
$ws.Range("T4").Value = 0.06246501555935044

$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("M5").Value = 0.3216333333333333
$ws.Range("N5").Value = 0.9649
$ws.Range("O5").Value = 0.1789992958012987
$ws.Range("P5").Value = 0.1789992958012987
$ws.Range("Q5").Value = 1.085454284366667
$ws.Range("R5").Value = 9.7690885593
$ws.Range("S5").Value = 0.01361898198157305
$ws.Range("T5").Value = 0.01361898198157305

$ws.Range("G6").Value = 4.037305666666668
$ws.Range("H6").Value = 12.111917
$ws.Range("I6").Value = 0.09101950487259411
$ws.Range("J6").Value = 0.09101950487259411
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.475208
$ws.Range("N6").Value = 4.425624
$ws.Range("O6").Value = 0.8210007041987012
$ws.Range("P6").Value = 0.8210007041987013
$ws.Range("Q6").Value = 5.955865617912002
$ws.Range("R6").Value = 53.60279056120801
$ws.Range("S6").Value = 0.07472707759621688
$ws.Range("T6").Value = 0.07472707759621688

$ws.Range("G7").Value = 4.037305666666668
$ws.Range("H7").Value = 12.111917
$ws.Range("I7").Value = 0.09101950487259411
$ws.Range("J7").Value = 0.09101950487259411
$ws.Range("M7").Value = 0.3216333333333333
$ws.Range("N7").Value = 0.9649
$ws.Range("O7").Value = 0.1789992958012987
$ws.Range("P7").Value = 0.1789992958012987
$ws.Range("Q7").Value = 1.298532079255556
$ws.Range("R7").Value = 11.6867887133
$ws.Range("S7").Value = 0.01629242727637722
$ws.Range("T7").Value = 0.01629242727637722
